$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.707.63'
$ws.Range("D3").Value = '3.255.98'
$ws.Range("E3").Value = '  +2.28%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.257.95'
$ws.Range("E8").Value = '  +2.43%  '
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("E11").Value = '  +4.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.508'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.44%  '
$ws.Range("D15").Value = '3.791.02'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '66.717.85'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '3.251.29'
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '507.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("E22").Value = '  +2.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("B27").Value = 'Hedera'
$ws.Range("C27").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.148'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +65.29%  '
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("E31").Value = '  -1.77%  '
$ws.Range("E32").Value = '  -7.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.16'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.13%  '
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '0.0₃0800'
$ws.Range("E37").Value = '  +16.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +19.25%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '55.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '495.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0429'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("E44").Value = '  -1.70%  '
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("D46").Value = '2.951.45'
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.49%  '
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.85%  '
